$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 4556.7
$ws.Range("D5").Value = 0.04754903036554232

$ws.Range("B6").Value = 8411.16
$ws.Range("D6").Value = 0.04293237925697658

$ws.Range("B7").Value = 8902.5
$ws.Range("D7").Value = 0.04867546569315735

$ws.Rows.Item(8).Delete()
